$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

# Remove the "bankdeposits" row (sheet row 9) from the variables table.
# This shifts every subsequent row up by one (table shrinks from 98 to 97 rows).
$ws.Range("A9").EntireRow.Delete()

# The conditional formatting ranges reference fixed row numbers and don't
# auto-shrink with the row delete, so re-point each one to its new extent.
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq '$K$1:$K$98') {
        $fc.ModifyAppliesToRange($ws.Range("K1:K97"))
    } elseif ($addr -eq '$O$94') {
        $fc.ModifyAppliesToRange($ws.Range("O93"))
    } elseif ($addr -eq '$H$2:$H$98') {
        $fc.ModifyAppliesToRange($ws.Range("H2:H97"))
    } elseif ($addr -eq '$O$2:$O$98') {
        $fc.ModifyAppliesToRange($ws.Range("O2:O97"))
    } elseif ($addr -eq '$P$2:$P$98') {
        $fc.ModifyAppliesToRange($ws.Range("P2:P97"))
    }
}

# Update the view: scroll back to the top and select D11 (matches the
# author's final cursor position after the edit).
$ws.Activate()
$ws.Range("D11").Select()
